$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

$ws.Range("A2").Value = -76.8833
$ws.Range("B2").Value = -76.807

$ws.Range("A3").Value = 38.0689
$ws.Range("B3").Value = 38.1275

$ws.Range("A4").Value = -75.6415
$ws.Range("B4").Value = -75.7172

$ws.Range("A5").Value = 39.5722
$ws.Range("B5").Value = 39.5136
